$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Refresh the cached text of the auto-updating "datetimeFigureOut" date
#    fields that live on the slide master and every slide layout (the value
#    PowerPoint caches the last time it rendered the field) from 5/6/2019 to
#    7/23/2019.
# ---------------------------------------------------------------------------
function Update-DateShapes($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq "5/6/2019") {
                $shp.TextFrame.TextRange.Text = "7/23/2019"
            }
        }
    }
}

Update-DateShapes $p.SlideMaster.Shapes
for ($li = 1; $li -le $p.SlideMaster.CustomLayouts.Count; $li++) {
    Update-DateShapes $p.SlideMaster.CustomLayouts.Item($li).Shapes
}

# ---------------------------------------------------------------------------
# 2) Remove the leftover screenshot pictures that were left on four slides
#    (installation/app screenshots no longer needed).
# ---------------------------------------------------------------------------
$picturesToRemove = @(
    @{slide = 2; id = 22},
    @{slide = 4; id = 17},
    @{slide = 5; id = 58},
    @{slide = 8; id = 18}
)

foreach ($target in $picturesToRemove) {
    $s = $p.Slides.Item($target.slide)
    for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
        $shp = $s.Shapes.Item($i)
        if ($shp.Id -eq $target.id) {
            $shp.Delete()
        }
    }
}
